# "Welcome and overview.pptx" update
#
# 1) Slide 3 ("We focus on the Wormhole") - content placeholder:
#    - add a new sub-bullet after "The next generation is the Blackhole"
#    - grow the placeholder to fit the extra line and turn on
#      "shrink text on overflow" (normAutofit)
# 2) Slide 5 ("Session plan") - timetable table: every session in the
#    afternoon now (shifted from a 9:00 start to a 14:00 start, i.e. +5h)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3: Wormhole / Blackhole bullet list
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$contentShape = $s3.Shapes.Item(2)

$tf = $contentShape.TextFrame
$tr = $tf.TextRange

$blackholePara = $tr.Paragraphs(4)
$blackholePara.InsertAfter("`rWe have both Wormhole and Blackhole, using Wormhole today") | Out-Null

# the newly inserted paragraph is now #5; indent it as a sub-bullet
# (IndentLevel is 1-based, so level 2 == <a:pPr lvl="1"/>)
$newPara = $tr.Paragraphs(5)
$newPara.IndentLevel = 2

# turn on shrink-text-on-overflow and resize the box to fit the new line
$tf.AutoSize = 2
$contentShape.Height = 223.3789

# ---------------------------------------------------------------------
# Slide 5: timetable, shift every slot by +5 hours (9:00 -> 14:00 ...)
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tableShape = $s5.Shapes.Item(2)
$tbl = $tableShape.Table

$dash = [char]0x2013

$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text  = "14:00 $dash 14:05"
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text  = "14:05 $dash 14:30"
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text  = "14:30 $dash 14:40"
$tbl.Cell(5, 1).Shape.TextFrame.TextRange.Text  = "14:40 $dash 15:30"
$tbl.Cell(6, 1).Shape.TextFrame.TextRange.Text  = "15:30 $dash 16:00"
$tbl.Cell(7, 1).Shape.TextFrame.TextRange.Text  = "16:00 $dash 16:05"
$tbl.Cell(8, 1).Shape.TextFrame.TextRange.Text  = "16:05 $dash 16:25"
$tbl.Cell(9, 1).Shape.TextFrame.TextRange.Text  = "16:25 $dash 17:25"
$tbl.Cell(10, 1).Shape.TextFrame.TextRange.Text = "17:25 $dash 17:30"
